# Adding the changes we made on may 9th
# Prepend 11 new rows of gyroscope samples above the existing data block
# (rows 2-12), which pushes the previous data down by 11 rows. Because the
# sheet's data window is capped at 30 rows, the oldest row (previously the
# last row of data) falls off the bottom and is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to insert at the top of the data (below the header row)
$newRows = @(
  @(0.007011067026972018, -0.03245915641838846, 0.03166781107641074),
  @(-0.0004164990228177805, -0.01167585594918232, -0.0006108652715655517),
  @(-0.0122173046693205, 0.0074830991216003, 0.0039706239476799),
  @(-0.04220523541285226, 0.007899598006836298, 0.001610462892461813),
  @(-0.01571589649062261, -0.009926560250195525, 0.002873843345283052),
  @(0.007719115523452058, -0.02122756669467137, 0.02797485068863757),
  @(-0.01150925694541499, -0.008524346549000439, 0.01102334066209471),
  @(-0.01870080676268442, -0.01660442801816813, -0.0122173046693205),
  @(-0.01634064570746632, -0.02958531457592136, -0.03082092817534093),
  @(-0.01731247691945591, -0.0009301814504645012, -0.02040845257314771),
  @(-0.03719336404041793, 0.1217704361135306, 0.03431951999664297)
)

$insertCount = $newRows.Length

# Shift the existing data (rows 2..21) down by inserting blank rows above it
$lastInsertRow = 1 + $insertCount
$ws.Range("A2:C$lastInsertRow").EntireRow.Insert()

# The insert copies formatting down from the header row; the source data
# rows carry no explicit styling, so clear it back off the new rows.
$ws.Range("A2:C$lastInsertRow").ClearFormats()

# Fill the newly-inserted rows with the new data
for ($i = 0; $i -lt $newRows.Length; $i++) {
  $row = $newRows[$i]
  $targetRow = $i + 2
  $ws.Cells.Item($targetRow, 1).Value = $row[0]
  $ws.Cells.Item($targetRow, 2).Value = $row[1]
  $ws.Cells.Item($targetRow, 3).Value = $row[2]
}

# The insert pushed the old last data row past the bottom of the kept
# window (it is now at row 32); remove that overflow row entirely.
$oldLastDataRow = 21
$overflowRow = $oldLastDataRow + $insertCount
$ws.Rows.Item($overflowRow).Delete()

Write-Output "done"
